$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1436.375
$ws.Range("I2").Value = 96.75
$ws.Range("J2").Value = 2776
$ws.Range("K2").Value = 96.75
$ws.Range("L2").Value = 2776
$ws.Range("M2").Value = 16.25
$ws.Range("N2").Value = -3002

$ws.Range("H6").Value = 22.5
$ws.Range("I6").Value = 30
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 90
$ws.Range("L6").Value = 45
$ws.Range("M6").Value = 22
$ws.Range("N6").Value = -269

$ws.Range("H9").Value = 1496.625
$ws.Range("I9").Value = 1478.3572
$ws.Range("J9").Value = 1624.5
$ws.Range("K9").Value = 1478.3572
$ws.Range("L9").Value = 1624.5
$ws.Range("M9").Value = -1309.3572
$ws.Range("N9").Value = -1962.5

$ws.Range("H12").Value = 500.75
$ws.Range("I12").Value = 500.75
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 500.75
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -330.75

$ws.Range("H29").Value = 7000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 7000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 21000
$ws.Range("N29").Value = -21562

$ws.Range("H38").Value = 158.2
$ws.Range("I38").Value = 158.2
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 474.6
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()

$ws.Range("H39").Value = 92.59999999999999
$ws.Range("I39").Value = 92.59999999999999
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 277.8
$ws.Range("L39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("N39").ClearContents()

$ws.Range("H42").Value = 1727.5
$ws.Range("I42").Value = 1803.3334
$ws.Range("J42").Value = 1500
$ws.Range("K42").Value = 5410.0002
$ws.Range("L42").Value = 4500
$ws.Range("M42").Value = -5180.0002
$ws.Range("N42").Value = -4960

$ws.Range("H51").Value = 17274.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 17274.75
$ws.Range("K51").Value = 0
$ws.Range("L51").ClearContents()
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -18242.75

$ws.Range("H54").Value = 6999.5
$ws.Range("I54").Value = 6999.5
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 6999.5
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

$ws.Range("H138").Value = 3748.0232
$ws.Range("I138").Value = 1911.3125
$ws.Range("J138").Value = 4167.843
$ws.Range("K138").Value = 5733.9375
$ws.Range("L138").Value = 12503.529
$ws.Range("M138").Value = -593.9375
$ws.Range("N138").Value = -22783.529

$ws.Range("H141").Value = 1939
$ws.Range("I141").Value = 1627.4286
$ws.Range("J141").Value = 2666
$ws.Range("K141").Value = 4882.2858
$ws.Range("L141").Value = 7998
$ws.Range("M141").Value = 297.7142000000003
$ws.Range("N141").Value = -18358

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1588925.5
$ws.Range("I32").Value = 805.59015
$ws.Range("J32").Value = 12352850
$ws.Range("K32").Value = 805.59015
$ws.Range("L32").Value = 12352850
$ws.Range("M32").Value = -518.59015
$ws.Range("N32").Value = -12353424

$ws.Range("H45").Value = 2182.9167
$ws.Range("I45").Value = 1819.6
$ws.Range("J45").Value = 3999.5
$ws.Range("K45").Value = 1819.6
$ws.Range("L45").Value = 3999.5
$ws.Range("M45").Value = -1442.6
$ws.Range("N45").Value = -4753.5

$ws.Range("H74").Value = 35717624
$ws.Range("I74").Value = 58826130
$ws.Range("J74").Value = 4480.1816
$ws.Range("K74").Value = 58826130
$ws.Range("L74").Value = 4480.1816
$ws.Range("M74").Value = -58825256
$ws.Range("N74").Value = -6228.1816

$ws.Range("H77").Value = 35717624
$ws.Range("I77").Value = 58826130
$ws.Range("J77").Value = 4480.1816
$ws.Range("K77").Value = 294130650
$ws.Range("L77").Value = 22400.908
$ws.Range("M77").Value = -294126282
$ws.Range("N77").Value = -31136.908

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 603987.8
$ws.Range("I132").Value = 682533.9
$ws.Range("J132").Value = 132711.28
$ws.Range("K132").Value = 2047601.7
$ws.Range("L132").Value = 398133.84
$ws.Range("M132").Value = -2045071.7
$ws.Range("N132").Value = -403193.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4993.857
$ws.Range("I31").Value = 1301.6364
$ws.Range("J31").Value = 9055.299999999999
$ws.Range("K31").Value = 1301.6364
$ws.Range("L31").Value = 9055.299999999999
$ws.Range("M31").Value = -1006.6364
$ws.Range("N31").Value = -9645.299999999999

$ws.Range("H34").Value = 4993.857
$ws.Range("I34").Value = 1301.6364
$ws.Range("J34").Value = 9055.299999999999
$ws.Range("K34").Value = 1301.6364
$ws.Range("L34").Value = 9055.299999999999
$ws.Range("M34").Value = -1099.6364
$ws.Range("N34").Value = -9459.299999999999

$ws.Range("H39").Value = 7699.3335
$ws.Range("I39").Value = 7699.3335
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 7699.3335
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -7308.3335

$ws.Range("H49").Value = 7699.3335
$ws.Range("I49").Value = 7699.3335
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 7699.3335
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -7517.3335

$ws.Range("H86").Value = 8342.795
$ws.Range("I86").Value = 6761.1055
$ws.Range("J86").Value = 9845.4
$ws.Range("K86").Value = 6761.1055
$ws.Range("L86").Value = 9845.4
$ws.Range("M86").Value = -5638.1055
$ws.Range("N86").Value = -12091.4

$ws.Range("H89").Value = 8342.795
$ws.Range("I89").Value = 6761.1055
$ws.Range("J89").Value = 9845.4
$ws.Range("K89").Value = 33805.5275
$ws.Range("L89").Value = 49227
$ws.Range("M89").Value = -28189.5275
$ws.Range("N89").Value = -60459

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 90.36364
$ws.Range("I38").Value = 20.75
$ws.Range("J38").Value = 276
$ws.Range("K38").Value = 62.25
$ws.Range("L38").Value = 828
$ws.Range("M38").Value = 284.75
$ws.Range("N38").Value = -1522

$ws.Range("H132").Value = 2499.111
$ws.Range("I132").Value = 1496.3334
$ws.Range("J132").Value = 2590.2727
$ws.Range("K132").Value = 13467.0006
$ws.Range("L132").Value = 23312.4543
$ws.Range("M132").Value = -10937.0006
$ws.Range("N132").Value = -28372.4543

$ws.Range("H133").Value = 4575
$ws.Range("I133").Value = 4575
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 13725
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -8665

$ws.Range("H137").Value = 4884.778
$ws.Range("I137").Value = 7950
$ws.Range("J137").Value = 2432.6
$ws.Range("K137").Value = 23850
$ws.Range("L137").Value = 7297.799999999999
$ws.Range("M137").Value = -18750
$ws.Range("N137").Value = -17497.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 23.266666
$ws.Range("I2").Value = 23.928572
$ws.Range("J2").Value = 14
$ws.Range("K2").Value = 23.928572
$ws.Range("L2").Value = 14
$ws.Range("M2").Value = 89.071428
$ws.Range("N2").Value = -240

$ws.Range("H70").Value = 8927.706
$ws.Range("I70").Value = 7761.3
$ws.Range("J70").Value = 10594
$ws.Range("K70").Value = 7761.3
$ws.Range("L70").Value = 10594
$ws.Range("M70").Value = -7491.3
$ws.Range("N70").Value = -11134

$ws.Range("H73").Value = 8927.706
$ws.Range("I73").Value = 7761.3
$ws.Range("J73").Value = 10594
$ws.Range("K73").Value = 7761.3
$ws.Range("L73").Value = 10594
$ws.Range("M73").Value = -6825.3
$ws.Range("N73").Value = -12466

$ws.Range("H80").Value = 11429.454
$ws.Range("I80").Value = 7100.6
$ws.Range("J80").Value = 15036.833
$ws.Range("K80").Value = 7100.6
$ws.Range("L80").Value = 15036.833
$ws.Range("M80").Value = -6102.6
$ws.Range("N80").Value = -17032.833

$ws.Range("H83").Value = 11429.454
$ws.Range("I83").Value = 7100.6
$ws.Range("J83").Value = 15036.833
$ws.Range("K83").Value = 35503
$ws.Range("L83").Value = 75184.16500000001
$ws.Range("M83").Value = -30511
$ws.Range("N83").Value = -85168.16500000001

$ws.Range("H97").Value = 912.9722
$ws.Range("I97").Value = 918.14813
$ws.Range("J97").Value = 897.44446
$ws.Range("K97").Value = 918.14813
$ws.Range("L97").Value = 897.44446
$ws.Range("M97").Value = -422.14813
$ws.Range("N97").Value = -1889.44446

$ws.Range("H132").Value = 142915500
$ws.Range("I132").Value = 333467230
$ws.Range("J132").Value = 1699.25
$ws.Range("K132").Value = 1000401690
$ws.Range("L132").Value = 5097.75
$ws.Range("M132").Value = -1000399160
$ws.Range("N132").Value = -10157.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3895.9285
$ws.Range("I122").Value = 3051.7144
$ws.Range("J122").Value = 6428.5713
$ws.Range("K122").Value = 9155.143199999999
$ws.Range("L122").Value = 19285.7139
$ws.Range("M122").Value = -6705.143199999999
$ws.Range("N122").Value = -24185.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 8345485.5
$ws.Range("I32").Value = 8345485.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 8345485.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -8345168.5

$ws.Range("H100").Value = 798.4286
$ws.Range("I100").Value = 497.2
$ws.Range("J100").Value = 1551.5
$ws.Range("K100").Value = 994.4
$ws.Range("L100").Value = 3103
$ws.Range("M100").Value = -453.4
$ws.Range("N100").Value = -4185
